$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 2000
$ws.Range("C3").Value = 3000
$ws.Range("C4").Value = 5000
